# Delete row 222 (CUR / Willemstad, Curacao) entirely.
# Excel will automatically shift all rows below it up by one,
# which matches the target diff (old row 223 "CAW" becomes new row 222, etc.,
# ending with old row 325 "SUV" becoming new row 324).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(222).Delete()
